$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a literal text value (e.g. "71.6%") into a cell without
# letting Excel auto-convert the percent-looking string into a numeric
# percentage (which would also swap in a brand new number-format style).
# We temporarily force a text format to keep the literal string, then
# restore the original look (General / centered) by pasting the format
# from a same-styled donor cell that is never otherwise touched.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($sheet, [string]$addr, [string]$val, [string]$donor)
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
    $sheet.Range($donor).Copy() | Out-Null
    $sheet.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) Class Statistics block (K/L columns) - recorded/missing counts + % values
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 197
$ws.Range("L7").Value = 35
Set-TextValue $ws "L9" "71.6%" "K9"
Set-TextValue $ws "L10" "75.1%" "K10"

# ---------------------------------------------------------------------------
# 2) Group Statistics block (M..S columns) for rows 16, 18, 19
# ---------------------------------------------------------------------------
$ws.Range("O16").Value = 43
$ws.Range("P16").Value = 4
Set-TextValue $ws "R16" "78.2%" "R15"
Set-TextValue $ws "S16" "77.3%" "S15"

$ws.Range("O18").Value = 42
$ws.Range("P18").Value = 6
Set-TextValue $ws "R18" "76.4%" "R15"
Set-TextValue $ws "S18" "78.1%" "S15"

$ws.Range("O19").Value = 40
$ws.Range("P19").Value = 7
Set-TextValue $ws "R19" "72.7%" "R15"
Set-TextValue $ws "S19" "66.6%" "S15"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) "Recorded By" (column G) - swap order of the two comma separated emails
# ---------------------------------------------------------------------------
$godaRows = @(32,33,34,37,38,39,87,89,90,92,93,94,95,96,197,198,199,201,202,203,205,206,253,255,258,260,261)
foreach ($r in $godaRows) {
    $ws.Range("G$r").Value = "160534@med.asu.edu.eg, emp17.mariam.m.goda@gmail.com"
}

$ws.Range("G125").Value = "ahmedali78112@gmail.com, abdallahashraf2023@gmail.com"

foreach ($r in @(187, 245)) {
    $ws.Range("G$r").Value = "emp17.nada.h.attia@gmail.com, dr.mohabelsawy@gmail.com"
}

# ---------------------------------------------------------------------------
# 4) Rows that flip from "Not Recorded" to "Recorded" (98, 173, 233)
#    - copy the green "Recorded" look from a neighbouring recorded row
#    - fill in Recorded By / Students / Status
# ---------------------------------------------------------------------------
$ws.Range("A97:I97").Copy() | Out-Null
$ws.Range("A98:I98").PasteSpecial(-4122) | Out-Null
$ws.Range("G98").Value = "emp17.mariam.a.saleh@gmail.com"
$ws.Range("H98").Value = "46/61"
$ws.Range("I98").Value = "Recorded"

$ws.Range("A172:I172").Copy() | Out-Null
$ws.Range("A173:I173").PasteSpecial(-4122) | Out-Null
$ws.Range("G173").Value = "awadayman129@gmail.com"
$ws.Range("H173").Value = "54/62"
$ws.Range("I173").Value = "Recorded"

$ws.Range("A232:I232").Copy() | Out-Null
$ws.Range("A233:I233").PasteSpecial(-4122) | Out-Null
$ws.Range("G233").Value = "abdallahashraf2023@gmail.com"
$ws.Range("H233").Value = "49/63"
$ws.Range("I233").Value = "Recorded"

$excel.CutCopyMode = 0
